# Add a new sheet "2025-03-24" at the end of the workbook (after the last existing sheet)
# mirroring the structure of the preceding daily price-summary sheets.
$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)

$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2025-03-24"

# Match the outline conventions used by the other daily sheets.
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

$ws.Cells.Item(1,1).Value = 'Match'
$ws.Cells.Item(1,2).Value = 'Seat Type'
$ws.Cells.Item(1,3).Value = 'Min_Price'
$ws.Cells.Item(1,4).Value = 'Avg_Price'
$ws.Cells.Item(1,5).Value = 'Ticket_Count'
$ws.Cells.Item(2,1).Value = 'Arsenal - Fulham'
$ws.Cells.Item(2,2).Value = 'Shortside Upper'
$ws.Cells.Item(2,3).Value = 0
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 57
$ws.Cells.Item(3,1).Value = 'Arsenal - Fulham'
$ws.Cells.Item(3,2).Value = 'Shortside Lower'
$ws.Cells.Item(3,3).Value = 0
$ws.Cells.Item(3,4).Value = 0
$ws.Cells.Item(3,5).Value = 63
$ws.Cells.Item(4,1).Value = 'Arsenal - Fulham'
$ws.Cells.Item(4,2).Value = 'Longside Upper'
$ws.Cells.Item(4,3).Value = 0
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,5).Value = 57
$ws.Cells.Item(5,1).Value = 'Arsenal - Fulham'
$ws.Cells.Item(5,2).Value = 'Longside Lower'
$ws.Cells.Item(5,3).Value = 0
$ws.Cells.Item(5,4).Value = 0
$ws.Cells.Item(5,5).Value = 113
$ws.Cells.Item(6,1).Value = 'Arsenal - Fulham'
$ws.Cells.Item(6,2).Value = 'Club Level'
$ws.Cells.Item(6,3).Value = 0
$ws.Cells.Item(6,4).Value = 0
$ws.Cells.Item(6,5).Value = 16
$ws.Cells.Item(7,1).Value = 'Arsenal - Fulham'
$ws.Cells.Item(7,2).Value = 'VIP & Executive Box'
$ws.Cells.Item(7,3).Value = 0
$ws.Cells.Item(7,4).Value = 0
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(8,1).Value = 'Arsenal - Brentford'
$ws.Cells.Item(8,2).Value = 'Longside Upper'
$ws.Cells.Item(8,3).Value = 0
$ws.Cells.Item(8,4).Value = 0
$ws.Cells.Item(8,5).Value = 57
$ws.Cells.Item(9,1).Value = 'Arsenal - Brentford'
$ws.Cells.Item(9,2).Value = 'Shortside Lower'
$ws.Cells.Item(9,3).Value = 0
$ws.Cells.Item(9,4).Value = 0
$ws.Cells.Item(9,5).Value = 57
$ws.Cells.Item(10,1).Value = 'Arsenal - Brentford'
$ws.Cells.Item(10,2).Value = 'Shortside Upper'
$ws.Cells.Item(10,3).Value = 0
$ws.Cells.Item(10,4).Value = 0
$ws.Cells.Item(10,5).Value = 50
$ws.Cells.Item(11,1).Value = 'Arsenal - Brentford'
$ws.Cells.Item(11,2).Value = 'Longside Lower'
$ws.Cells.Item(11,3).Value = 0
$ws.Cells.Item(11,4).Value = 0
$ws.Cells.Item(11,5).Value = 77
$ws.Cells.Item(12,1).Value = 'Arsenal - Brentford'
$ws.Cells.Item(12,2).Value = 'Away Section'
$ws.Cells.Item(12,3).Value = 0
$ws.Cells.Item(12,4).Value = 0
$ws.Cells.Item(12,5).Value = 1
$ws.Cells.Item(13,1).Value = 'Arsenal - Brentford'
$ws.Cells.Item(13,2).Value = 'Club Level'
$ws.Cells.Item(13,3).Value = 0
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 18
$ws.Cells.Item(14,1).Value = 'Arsenal - Crystal Palace'
$ws.Cells.Item(14,2).Value = 'Shortside Upper'
$ws.Cells.Item(14,3).Value = 0
$ws.Cells.Item(14,4).Value = 0
$ws.Cells.Item(14,5).Value = 87
$ws.Cells.Item(15,1).Value = 'Arsenal - Crystal Palace'
$ws.Cells.Item(15,2).Value = 'Longside Lower'
$ws.Cells.Item(15,3).Value = 0
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = 111
$ws.Cells.Item(16,1).Value = 'Arsenal - Crystal Palace'
$ws.Cells.Item(16,2).Value = 'Shortside Lower'
$ws.Cells.Item(16,3).Value = 0
$ws.Cells.Item(16,4).Value = 0
$ws.Cells.Item(16,5).Value = 90
$ws.Cells.Item(17,1).Value = 'Arsenal - Crystal Palace'
$ws.Cells.Item(17,2).Value = 'Longside Upper'
$ws.Cells.Item(17,3).Value = 0
$ws.Cells.Item(17,4).Value = 0
$ws.Cells.Item(17,5).Value = 82
$ws.Cells.Item(18,1).Value = 'Arsenal - Crystal Palace'
$ws.Cells.Item(18,2).Value = 'Away Section'
$ws.Cells.Item(18,3).Value = 0
$ws.Cells.Item(18,4).Value = 0
$ws.Cells.Item(18,5).Value = 1
$ws.Cells.Item(19,1).Value = 'Arsenal - Crystal Palace'
$ws.Cells.Item(19,2).Value = 'Club Level'
$ws.Cells.Item(19,3).Value = 0
$ws.Cells.Item(19,4).Value = 0
$ws.Cells.Item(19,5).Value = 19
$ws.Cells.Item(20,1).Value = 'Arsenal v Real Madrid : Champions League 2024-2025'
$ws.Cells.Item(20,2).Value = 'Shortside Upper'
$ws.Cells.Item(20,3).Value = 0
$ws.Cells.Item(20,4).Value = 0
$ws.Cells.Item(20,5).Value = 38
$ws.Cells.Item(21,1).Value = 'Arsenal v Real Madrid : Champions League 2024-2025'
$ws.Cells.Item(21,2).Value = 'Longside Upper'
$ws.Cells.Item(21,3).Value = 0
$ws.Cells.Item(21,4).Value = 0
$ws.Cells.Item(21,5).Value = 60
$ws.Cells.Item(22,1).Value = 'Arsenal v Real Madrid : Champions League 2024-2025'
$ws.Cells.Item(22,2).Value = 'Shortside Lower'
$ws.Cells.Item(22,3).Value = 0
$ws.Cells.Item(22,4).Value = 0
$ws.Cells.Item(22,5).Value = 36
$ws.Cells.Item(23,1).Value = 'Arsenal v Real Madrid : Champions League 2024-2025'
$ws.Cells.Item(23,2).Value = 'Longside Lower'
$ws.Cells.Item(23,3).Value = 0
$ws.Cells.Item(23,4).Value = 0
$ws.Cells.Item(23,5).Value = 92
$ws.Cells.Item(24,1).Value = 'Arsenal v Real Madrid : Champions League 2024-2025'
$ws.Cells.Item(24,2).Value = 'Club Level'
$ws.Cells.Item(24,3).Value = 0
$ws.Cells.Item(24,4).Value = 0
$ws.Cells.Item(24,5).Value = 30
$ws.Cells.Item(25,1).Value = 'Arsenal v Real Madrid : Champions League 2024-2025'
$ws.Cells.Item(25,2).Value = 'Away Section'
$ws.Cells.Item(25,3).Value = 0
$ws.Cells.Item(25,4).Value = 0
$ws.Cells.Item(25,5).Value = 4
$ws.Cells.Item(26,1).Value = 'Arsenal v Real Madrid : Champions League 2024-2025'
$ws.Cells.Item(26,2).Value = 'VIP & Executive Box'
$ws.Cells.Item(26,3).Value = 0
$ws.Cells.Item(26,4).Value = 0
$ws.Cells.Item(26,5).Value = 2

# Header formatting to match the other daily sheets (bold, centered, thin boxed border)
$header = $ws.Range("A1:E1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

$ws.Range("A1").Select()

